$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Tue Nov 12 18:02:58 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 18:03:12 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 18:03:25 EST 2024"
$ws.Range("B5").Value = "Tue Nov 12 18:03:38 EST 2024"
